# Add creditsGiven field to game models and implement credit card assignment
# script -- applied here as the corresponding spreadsheet edits:
#   1) Expand the Powerups list already granted to the "IMPOSTORS" team
#      (row 24, column G) with the newly-assigned powerups 8, 9, 7, 6.
#   2) Append a new team "Reapers" (with its roster) as a new data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the Powerups cell for the IMPOSTORS team (row 24).
$ws.Range("G24").Value = "1, 2, 3, 4, 5, 8, 9, 7, 6"

# 2) Append the new "Reapers" team on the next free row.
$newRow = $ws.Range("A" + $ws.Rows.Count()).End(-4162).Row() + 1

$ws.Cells.Item($newRow, 1).Value = "Reapers"
$ws.Cells.Item($newRow, 2).Value = "Ishwinder, Aradhya, Kumud, Arnav"
$ws.Cells.Item($newRow, 3).Value = ""
$ws.Cells.Item($newRow, 4).Value = ""
$ws.Cells.Item($newRow, 5).Value = ""
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = ""
